$d = $word.ActiveDocument

# --- 0.1 Course overview heading: Heading2 -> Heading3, renumber 0.1 -> 0.0.1 ---
$pOverview = $d.Paragraphs(7)
$pOverview.Style = "Heading 3"
$tab = "`t"
$pOverview.Range.Find.Execute("0.1$tab", $true, $false, $false, $false, $false, $false, 1, $false, "0.0.1$tab", 2) | Out-Null

# --- Week 1 block ---
$d.Paragraphs(8).Range.Find.Execute("Week 1", $true, $false, $false, $false, $false, $true, 1, $false, "Week 1 Ethics and Law", 2) | Out-Null
$d.Paragraphs(9).Range.Find.Execute("Ethics and Law – ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Week 2 block ---
$d.Paragraphs(10).Range.Find.Execute("Week 2", $true, $false, $false, $false, $false, $true, 1, $false, "Week 2 Crime and Justice", 2) | Out-Null
$d.Paragraphs(11).Range.Find.Execute("Crime and Justice – ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Week 3 block ---
$d.Paragraphs(12).Range.Find.Execute("Week 3", $true, $false, $false, $false, $false, $true, 1, $false, "Week 3 Home and City", 2) | Out-Null
$d.Paragraphs(13).Range.Find.Execute("Home and City – ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Week 4 block ---
$d.Paragraphs(14).Range.Find.Execute("Week 4", $true, $false, $false, $false, $false, $true, 1, $false, "Week 4 Money and Markets", 2) | Out-Null
$d.Paragraphs(15).Range.Find.Execute("Money and Markets – ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Week 5 block ---
$d.Paragraphs(16).Range.Find.Execute("Week 5", $true, $false, $false, $false, $false, $true, 1, $false, "Week 5 Life and Health", 2) | Out-Null
$d.Paragraphs(17).Range.Find.Execute("Life and Health – ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- Last updated timestamp ---
$d.Paragraphs(18).Range.Find.Execute("15:21", $true, $false, $false, $false, $false, $true, 1, $false, "15:40", 2) | Out-Null
